# Dataset.xlsx cleanup: remove the two duplicate/blank "TestCase10" rows and
# re-label the now-contiguous "TestCase07" / "TestCase08" data rows so each
# test case has exactly 3 parameter rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 26:27 held the duplicate "TestCase10" / "[]" placeholder rows - drop them.
# (Everything below shifts up by two rows as a result.)
$ws.Rows("26:27").Delete()

# Column A ("Test Case Name") for the six rows that used to belong to the stray
# "TestCase09" / "TestCase10" groups now needs to read "TestCase07" (x3) and
# "TestCase08" (x3) so the labels line up with their parameter rows again.
$ws.Range("A20").Value = "TestCase07"
$ws.Range("A21").Value = "TestCase07"
$ws.Range("A22").Value = "TestCase07"
$ws.Range("A23").Value = "TestCase08"
$ws.Range("A24").Value = "TestCase08"
$ws.Range("A25").Value = "TestCase08"

# Rows 20, 22, 23, 24 and 25 pick up the same cell formatting already used by
# the rows below them (e.g. A29) - copy it across so the look matches.
$fmtSrc = $ws.Range("A29")
$fmtSrc.Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("A24").PasteSpecial(-4122)
$ws.Range("A25").PasteSpecial(-4122)
